$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 80, shifting existing rows 80-99 down to 81-100
$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with the new weekly record
$ws.Range("A80").Value = 6
$ws.Range("B80").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C80").Value = "Metropolitana"
$ws.Range("D80").Value = 44809
$ws.Range("D80").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E80").Value = 13
$ws.Range("F80").Value = 100114007
$ws.Range("G80").Value = "Jengibre"
$ws.Range("H80").Value = "Sin especificar"
$ws.Range("I80").Value = "Primera"
$ws.Range("J80").Value = 280
$ws.Range("K80").Value = 10000
$ws.Range("L80").Value = 11000
$ws.Range("M80").Value = 10571
$ws.Range("N80").Value = "$/caja 13 kilos"
$ws.Range("O80").Value = "Perú"
$ws.Range("P80").Value = 813
$ws.Range("Q80").Value = 13
$ws.Range("R80").Value = "Hortaliza"
